$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Report" to "Sheet1"
$ws.Name = "Sheet1"

# Update membership-count column (B) values for the society rows
$ws.Range("B2").Value = 275
$ws.Range("B3").Value = 250
$ws.Range("B4").Value = 184
$ws.Range("B5").Value = 85
$ws.Range("B6").Value = 400
